$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.460.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.50%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07730"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7386"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.873.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.132"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "272.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.441.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007592"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.124.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.250"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.191"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.337"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("E28").Value = "  +2.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.371"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09950"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +4.75%  "
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.074"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04783"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7013"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.715"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.731"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.347"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4204"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8385"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.302"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.090"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "928.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05644"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
